$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.305.03"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.390.01"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.00"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.47"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("E9").Value = "  +7.60%  "
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("E11").Value = "  +3.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000283"
$ws.Range("E12").Value = "  +4.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "679.79"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.935.64"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.406.55"
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.392.72"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.73"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.19"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.20"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.67"
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.90"
$ws.Range("E28").Value = "  +2.53%  "
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "556.00"
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.60"
$ws.Range("E33").Value = "  +9.74%  "
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.58"
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.669.10"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  +5.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.57"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0723"
$ws.Range("E40").Value = "  +7.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.27"
$ws.Range("E41").Value = "  +3.47%  "
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0424"
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.30"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.69"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("E48").Value = "  +5.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.04"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.69"
$ws.Range("E51").Value = "  +5.53%  "
